$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: rename the "테스트1" label to "테스트" and bump the maturity date
# (H17) from 2025-08-01 (45870) to 2025-08-29 (45898).
$ws.Range("B17").Value = "테스트"
$ws.Range("H17").Value = 45898

# Row 18 (the old "테스트2"/A05 entry) is removed entirely.
$ws.Rows(18).Delete()

# Reflect the author's last selection in the saved view.
[void]$ws.Range("H18").Select()
